# Update "想去人数" (want-to-go count) values in column F across sheets
# as generated by the gh-pages data refresh (commit 456a3b4).
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1202
$ws.Range("F4").Value = 12687
$ws.Range("F5").Value = 714
$ws.Range("F6").Value = 59
$ws.Range("F10").Value = 1857
$ws.Range("F11").Value = 39
$ws.Range("F13").Value = 510
$ws.Range("F16").Value = 338
$ws.Range("F17").Value = 225
$ws.Range("F18").Value = 287
$ws.Range("F19").Value = 126
$ws.Range("F20").Value = 123
$ws.Range("F21").Value = 27
$ws.Range("F22").Value = 214
$ws.Range("F23").Value = 239
$ws.Range("F24").Value = 1270
$ws.Range("F25").Value = 329
$ws.Range("F26").Value = 63

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 279
$ws.Range("F5").Value = 4450
$ws.Range("F6").Value = 142
$ws.Range("F16").Value = 11

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 865
$ws.Range("F3").Value = 2008

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 865
$ws.Range("F6").Value = 1203
$ws.Range("F7").Value = 12687
$ws.Range("F8").Value = 279
$ws.Range("F9").Value = 714
$ws.Range("F10").Value = 2009
$ws.Range("F11").Value = 59
$ws.Range("F15").Value = 1857
$ws.Range("F16").Value = 39
$ws.Range("F18").Value = 510
$ws.Range("F19").Value = 4450
$ws.Range("F21").Value = 142
$ws.Range("F22").Value = 142
$ws.Range("F29").Value = 338
$ws.Range("F31").Value = 225
$ws.Range("F32").Value = 287
$ws.Range("F33").Value = 126
$ws.Range("F34").Value = 123
$ws.Range("F35").Value = 27
$ws.Range("F37").Value = 214
$ws.Range("F40").Value = 239
$ws.Range("F41").Value = 1270
$ws.Range("F42").Value = 11
$ws.Range("F43").Value = 329
$ws.Range("F44").Value = 63
